$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("HPtFM")

# --- "HPtFM" data sheet: replace "thermochemical water splitting" pathway
#     with "hydrocarbon partial oxidation" and mark it as using "heavy or
#     residual fuel oil" (column I) instead of no fuel at all.
$wsData.Range("A6").Value = "hydrocarbon partial oxidation"
$wsData.Range("I6").Value = 1

# --- "About" sheet: remove the now-obsolete explanatory note about
#     thermochemical water splitting (previously rows 14-16).
$wsAbout.Rows.Item(14).Resize(3).Delete() | Out-Null

# --- Restore the selections left in each sheet by the editing session.
$wsData.Range("A3").Select() | Out-Null
$wsAbout.Range("B9").Select() | Out-Null
